# Before discussion with gls
# Update the I-column "percent difference" formulas to express the
# result as a percentage (multiply by 100), add a new K2 summary
# formula (average of the H column), and leave the new selection
# on the recomputed I column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# I2 was its own (non-shared) formula in the original workbook -
# update it on its own so it stays a standalone formula cell.
$ws.Range("I2").Formula = "=H2/B2*100"

# I3:I87 were two shared-formula blocks (I3:I66 and I67:I87) - writing
# the same formula text across the whole range preserves that shared
# grouping while updating the underlying formula.
$ws.Range("I3:I87").Formula = "=H3/B3*100"

# New column K: average of the H column (note the range intentionally
# spills one row past the last data row, H88, matching the source edit).
$ws.Range("K2").Formula = "=AVERAGE(H2:H88)"

# Reflect the author's new selection (I2 active cell, I2:I87 selected).
$null = $ws.Range("I2:I87").Select()
